# Questions can be imported from excel sheet
#
# This applies the edits captured in the target diff:
#   1. Cell H2 on the "Input" sheet changes from the number 3 to the text
#      "3,1" (a new shared string is minted for it automatically).
#   2. The active selection moves from N18 to M17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# 1. correct_options for the first sample row now allows multiple answers
#    ("3,1") instead of a single numeric option ("3").
$ws.Range("H2").Value = "3,1"

# 2. Move/save the selection like the author's Excel session did.
$ws.Range("M17").Select()
